$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notas")

$v0 = @'
4 - Por determinação judicial, os municípios de Livramento e Taperoá, do estado da Paraíba, tiveram seus limites territoriais alterados passando a vigorar, para efeitos de distribuição do Fundo de Participação dos Municípios, as populações de 7.248 e 15.316 habitantes, respectivamente. Processo Judicial nº 0000301-96.2006.4.05.8200 - Tribunal de Justiça da Paraíba.

Em 2017:
'@
$ws.Range("A21").Value = $v0

$v1 = @'
1 - População judicial do município de Porto Velho-RO: 494.013 habitantes. Processo Judicial nº 12316-40.2016.4.01.4100 - Seção Judiciária de Rondônia.
'@
$ws.Range("A22").Value = $v1

$v2 = @'
2 - População judicial do município de Santa Isabel do Rio Negro-AM: entre 23.773 e 30.564 habitantes. Parecer de Força Executória nº 00007/2017/NUCOB-GEAC/PFAM/PGF/AGU, em trâmite na 3ª VF/AM.
'@
$ws.Range("A23").Value = $v2

$v3 = @'
3 - População judicial do município de Urucará-AM: entre 16.981 e 23.772 habitantes. Parecer de Força Executória nº 00004/2017/NUCOB-GEAC/PFAM/PGF/AGU, em trâmite na 3ª VF/AM.
'@
$ws.Range("A24").Value = $v3

$v4 = @'
4 - População judicial do município de Jacareacanga-PA: 41.487 habitantes. Processo Judicial nº 798-41.2011.4.01.3902, Seção Judiciária de Itaituba-PA.
'@
$ws.Range("A25").Value = $v4

$v5 = @'
5 - População judicial do município de Paço do Lumiar - MA: superior a 156.216 habitantes. Processo Judicial nº13916-98.2017.4.01.3700 - Seção Judiciária do Maranhão- MA.
'@
$ws.Range("A26").Value = $v5

$v6 = @'
6 - População judicial do município Livramento-PB: 7.262 habitantes. Processo Judicial nº 0000301-96.2006.4.05.8200 - Tribunal de Justiça da Paraíba - PB.
'@
$ws.Range("A27").Value = $v6

$v7 = @'
7 - População judicial do município de Taperoá-PB: 15.400 habitantes. Processo Judicial nº 0000301-96.2006.4.05.8200 - Tribunal de Justiça da Paraíba - PB.
'@
$ws.Range("A28").Value = $v7

$v8 = @'
8 - População judicial do município Coronel João Sá-BA: 17.422 habitantes. Processo Judicial nº 0002222-53.2017.4.01.3306 - Vara Única de Paulo Afonso-BA.
'@
$ws.Range("A29").Value = $v8

$v9 = @'
31 - População judicial do Município de Ribeirão do Pinhal-PR: 13.601 habitantes. Processo Judicial nº 5001464-05.2018.4.04.7013 – Tribunal Regional Federal da 4ª Região.

De 2018 a 2021:
As diferenças entre as populações das Unidades da Federação obtidas da soma das estimativas municipais e aquelas projetadas nas Projeções de População, Brasil e Unidades da Federação, Revisão 2018, devem-se à alteração de limites territoriais ocorridas entre os estados, após o Censo Demográfico 2010. 

 
Para o histórico de alterações, consulte o link https://www.ibge.gov.br/estatisticas-novoportal/sociais/populacao/9103-estimativas-de-populacao.html.
'@
$ws.Range("A106").Value = $v9

# Clear the stray empty-string cell so it becomes a genuinely blank cell
$ws.Range("A107").ClearContents()

# Materialize the blank cells that become part of the used range after a resave
$ws.Range("B1:B108").Style = "Normal"

$wsTabela = $wb.Worksheets.Item("Tabela")
$wsTabela.Range("A1:V2").Style = "Normal"
$wsTabela.Range("A3:V4").Style = "Normal"
$wsTabela.Range("A32:V32").Style = "Normal"
